$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "christopher nolan"
$ws.Range("A3").Value = "taika waititi"
$ws.Range("A4").Value = "james cameron"
$ws.Range("A5").Value = "steven speilberg"

$ws.Range("A2").RowHeight = 15

$ws.Range("A3").Select()
